$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 45 new rows (week ending 2021-04-18), columns: Datum, Receipt Number, Konto, Beskrivning, Debet, Kredit
$newRows = @(
    @(44298, "Reko80", 3011, "Reko Swish +46723047499", $null, 423.21),
    @(44298, "Reko80", 2611, "Reko Swish +46723047499", $null, 50.79),
    @(44298, "Reko80", 1930, "Reko Swish +46723047499", 474, $null),
    @(44298, "Reko81", 3011, "Reko Swish +46734058272", $null, 141.07),
    @(44298, "Reko81", 2611, "Reko Swish +46734058272", $null, 16.93),
    @(44298, "Reko81", 1930, "Reko Swish +46734058272", 158, $null),
    @(44298, "Reko82", 3011, "Reko Swish +46767736486", $null, 512.5),
    @(44298, "Reko82", 2611, "Reko Swish +46767736486", $null, 61.5),
    @(44298, "Reko82", 1930, "Reko Swish +46767736486", 574, $null),
    @(44298, "Reko83", 3011, "Reko Swish +46733706356", $null, 141.07),
    @(44298, "Reko83", 3011, "Reko Swish +46733706356", $null, 141.07),
    @(44298, "Reko83", 2611, "Reko Swish +46733706356", $null, 16.93),
    @(44298, "Reko84", 2611, "Reko Swish +46733706356", $null, 16.93),
    @(44298, "Reko84", 1930, "Reko Swish +46733706356", 158, $null),
    @(44298, "Reko84", 1930, "Reko Swish +46733706356", 158, $null),
    @(44299, "Reko85", 3011, "Reko Swish +46720444719", $null, 141.07),
    @(44299, "Reko85", 2611, "Reko Swish +46720444719", $null, 16.93),
    @(44299, "Reko85", 1930, "Reko Swish +46720444719", 158, $null),
    @(44299, "Reko86", 3011, "Reko Swish +46702331968", $null, 616.07),
    @(44299, "Reko86", 2611, "Reko Swish +46702331968", $null, 73.93),
    @(44299, "Reko86", 1930, "Reko Swish +46702331968", 690, $null),
    @(44299, "Reko87", 3011, "Reko Swish +46703564388", $null, 398.21),
    @(44299, "Reko87", 2611, "Reko Swish +46703564388", $null, 47.79),
    @(44299, "Reko87", 1930, "Reko Swish +46703564388", 446, $null),
    @(44299, "Reko88", 3011, "Reko Swish +46705716511", $null, 282.14),
    @(44299, "Reko88", 2611, "Reko Swish +46705716511", $null, 33.86),
    @(44299, "Reko88", 1930, "Reko Swish +46705716511", 316, $null),
    @(44299, "Reko89", 3011, "Reko Swish +46705725419", $null, 230.36),
    @(44299, "Reko89", 2611, "Reko Swish +46705725419", $null, 27.64),
    @(44299, "Reko89", 1930, "Reko Swish +46705725419", 258, $null),
    @(44300, "Reko90", 3011, "Reko Swish +46703384055", $null, 610.71),
    @(44300, "Reko90", 2611, "Reko Swish +46703384055", $null, 73.29),
    @(44300, "Reko90", 1930, "Reko Swish +46703384055", 684, $null),
    @(44301, "Reko91", 3011, "Reko Swish +46706649892", $null, 423.21),
    @(44301, "Reko91", 2611, "Reko Swish +46706649892", $null, 50.79),
    @(44301, "Reko91", 1930, "Reko Swish +46706649892", 474, $null),
    @(44302, "Reko92", 3011, "Reko Swish +46739582203", $null, 575.89),
    @(44302, "Reko92", 2611, "Reko Swish +46739582203", $null, 69.11),
    @(44302, "Reko92", 1930, "Reko Swish +46739582203", 645, $null),
    @(44302, $null, 5460, "IKEA BARKARBY K0135", 3104, $null),
    @(44302, $null, 2641, "IKEA BARKARBY K0135", 776, $null),
    @(44302, $null, 1930, "IKEA BARKARBY K0135", $null, 3880),
    @(44303, $null, 5460, "IKEA BARKARBY K0135", 14845.6, $null),
    @(44303, $null, 2641, "IKEA BARKARBY K0135", 3711.4, $null),
    @(44303, $null, 1930, "IKEA BARKARBY K0135", $null, 18557)
)

$dateFormat = $ws.Range("A604").NumberFormat
$startRow = 605
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Column A: Datum (date-formatted serial number, matches style used by prior rows)
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 1).Value = $data[0]

    # Column B: Receipt Number (text, blank for IKEA rows)
    if ($data[1] -eq $null) {
        $ws.Cells.Item($r, 2).Value = ""
    } else {
        $ws.Cells.Item($r, 2).Value = $data[1]
    }

    # Column C: Konto
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Column D: Beskrivning
    $ws.Cells.Item($r, 4).Value = $data[3]

    # Column E: Debet (blank when value not given)
    if ($data[4] -eq $null) {
        $ws.Cells.Item($r, 5).Value = ""
    } else {
        $ws.Cells.Item($r, 5).Value = $data[4]
    }

    # Column F: Kredit (blank when value not given)
    if ($data[5] -eq $null) {
        $ws.Cells.Item($r, 6).Value = ""
    } else {
        $ws.Cells.Item($r, 6).Value = $data[5]
    }
}

